$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk evaluation")

# Update risk probability values in column C; dependent formulas in column E
# (E = C*D, shared formula) recalc automatically.
$ws.Range("C7").Value = 0.18
$ws.Range("C8").Value = 0.2
$ws.Range("C9").Value = 0.1
$ws.Range("C10").Value = 0.18

$ws.Range("C12").Value = 0.005
$ws.Range("D12").Value = 2

# Update the view position / active selection to match the author's final
# scroll & selection state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select()
